$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44326
$ws.Range("M2").Value = 65
$ws.Range("D3").Value = 44326
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 67
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range("S3").Value = 800
$ws.Range("D5").Value = 44323
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 60
$ws.Range("D6").Value = 44323
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 50
$ws.Range("D7").Value = 44315
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 45
$ws.Range("N7").Value = 10000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 10000
$ws.Range("S7").Value = 1000
$ws.Range("D8").Value = 44333
$ws.Range("M8").Value = 58
$ws.Range("R8").Value = "Provincia de Quillota"
$ws.Range("D9").Value = 44333
$ws.Range("M9").Value = 65
$ws.Range("R9").Value = "Provincia de Quillota"
$ws.Range("D10").Value = 44333
$ws.Range("M10").Value = 60
$ws.Range("R10").Value = "Provincia de Quillota"
$ws.Range("D11").Value = 44321
$ws.Range("M11").Value = 58
$ws.Range("N11").Value = 9000
$ws.Range("O11").Value = 9000
$ws.Range("P11").Value = 9000
$ws.Range("S11").Value = 900
$ws.Range("D12").Value = 44307
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 40
$ws.Range("N12").Value = 10000
$ws.Range("O12").Value = 10000
$ws.Range("P12").Value = 10000
$ws.Range("S12").Value = 1000
$ws.Range("D13").Value = 44302
$ws.Range("M13").Value = 45
$ws.Range("D14").Value = 44309
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 45
$ws.Range("N14").Value = 10000
$ws.Range("O14").Value = 10000
$ws.Range("P14").Value = 10000
$ws.Range("S14").Value = 1000
$ws.Range("D15").Value = 44319
$ws.Range("M15").Value = 68
$ws.Range("D16").Value = 44319
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 57
$ws.Range("N16").Value = 8000
$ws.Range("O16").Value = 8000
$ws.Range("P16").Value = 8000
$ws.Range("S16").Value = 800
$ws.Range("D17").Value = 44308
$ws.Range("M17").Value = 45
$ws.Range("D18").Value = 44308
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 48
$ws.Range("N18").Value = 8000
$ws.Range("O18").Value = 8000
$ws.Range("P18").Value = 8000
$ws.Range("S18").Value = 800
$ws.Range("D19").Value = 44306
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 45
$ws.Range("N19").Value = 10000
$ws.Range("O19").Value = 10000
$ws.Range("P19").Value = 10000
$ws.Range("S19").Value = 1000
$ws.Range("D20").Value = 44343
$ws.Range("L20").Value = "Especial"
$ws.Range("M20").Value = 47
$ws.Range("R20").Value = "Región Metropolitana"
$ws.Range("D21").Value = 44343
$ws.Range("M21").Value = 50
$ws.Range("N21").Value = 9000
$ws.Range("O21").Value = 9000
$ws.Range("P21").Value = 9000
$ws.Range("R21").Value = "Región Metropolitana"
$ws.Range("S21").Value = 900
$ws.Range("D22").Value = 44343
$ws.Range("M22").Value = 58
$ws.Range("N22").Value = 8000
$ws.Range("O22").Value = 8000
$ws.Range("P22").Value = 8000
$ws.Range("R22").Value = "Región Metropolitana"
$ws.Range("S22").Value = 800
$ws.Range("D23").Value = 44329
$ws.Range("N23").Value = 9000
$ws.Range("O23").Value = 9000
$ws.Range("P23").Value = 9000
$ws.Range("R23").Value = "Región Metropolitana"
$ws.Range("S23").Value = 900
$ws.Range("D24").Value = 44329
$ws.Range("M24").Value = 50
$ws.Range("R24").Value = "Región Metropolitana"
$ws.Range("D25").Value = 44328
$ws.Range("M25").Value = 45
$ws.Range("N25").Value = 8000
$ws.Range("O25").Value = 8000
$ws.Range("P25").Value = 8000
$ws.Range("S25").Value = 800
$ws.Range("D26").Value = 44328
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 48
$ws.Range("N26").Value = 7000
$ws.Range("O26").Value = 7000
$ws.Range("P26").Value = 7000
$ws.Range("S26").Value = 700
$ws.Range("D27").Value = 44301
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 45
$ws.Range("N27").Value = 10000
$ws.Range("O27").Value = 10000
$ws.Range("P27").Value = 10000
$ws.Range("S27").Value = 1000
$ws.Range("D28").Value = 44322
$ws.Range("N28").Value = 10000
$ws.Range("O28").Value = 10000
$ws.Range("P28").Value = 10000
$ws.Range("R28").Value = "Provincia de Quillota"
$ws.Range("S28").Value = 1000
$ws.Range("D29").Value = 44322
$ws.Range("M29").Value = 40
$ws.Range("R29").Value = "Provincia de Quillota"
$ws.Range("D30").Value = 44312
$ws.Range("M30").Value = 48
